# Rename the "Trend (CAGR)" / "Trend (Reg)" forecast-method labels to the
# shorter "CAGR" / "Trend" labels used for week 12's content, across every
# column group on the "Sales Method" (row 4) and "COGS Method" (row 5)
# header rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - "Sales Method" column headers (one CAGR/Trend pair per forecast group)
$ws.Range("H4").Value = "CAGR"
$ws.Range("I4").Value = "Trend"
$ws.Range("L4").Value = "CAGR"
$ws.Range("M4").Value = "Trend"
$ws.Range("P4").Value = "CAGR"
$ws.Range("Q4").Value = "Trend"
$ws.Range("T4").Value = "CAGR"
$ws.Range("U4").Value = "Trend"

# Row 5 - "COGS Method" column headers
$ws.Range("H5").Value = "CAGR"
$ws.Range("I5").Value = "Trend"
$ws.Range("T5").Value = "CAGR"
$ws.Range("U5").Value = "CAGR"
$ws.Range("V5").Value = "CAGR"
$ws.Range("W5").Value = "CAGR"

# Drop the stale duplicate chart-tracking defined name left over from a
# previous save (both pointed at the same range; only one is still used by
# the chartEx plot)
$wb.Names.Item("_xlchart.v1.1").Delete()

# Restore the cursor to where the author left it
$ws.Range("A2").Select()
